# Update for 15 April: append a new "4/14/20" deaths column (AG) to the
# US states deaths table, mirroring the formatting of the previous day's
# column (AF).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the added date column.
$ws.Range("AG1").Value = " 4/14/20"

# New daily cumulative-deaths figures, keyed by worksheet row number.
# Row order follows the existing state ordering already present in column A.
$values = [ordered]@{
    2  = 114    # Alabama
    3  = 9      # Alaska
    4  = 131    # Arizona
    5  = 32     # Arkansas
    6  = 789    # California
    7  = 329    # Colorado
    8  = 671    # Connecticut
    9  = 43     # Delaware
    10 = 67     # District of Columbia
    11 = 571    # Florida
    12 = 524    # Georgia
    13 = 5      # Guam
    14 = 9      # Hawaii
    15 = 39     # Idaho
    16 = 868    # Illinois
    17 = 387    # Indiana
    18 = 49     # Iowa
    19 = 69     # Kansas
    20 = 115    # Kentucky
    21 = 1013   # Louisiana
    22 = 20     # Maine
    23 = 302    # Maryland
    24 = 957    # Massachusetts
    25 = 1768   # Michigan
    26 = 79     # Minnesota
    27 = 111    # Mississippi
    28 = 142    # Missouri
    29 = 7      # Montana
    30 = 20     # Nebraska
    31 = 130    # Nevada
    32 = 27     # New Hampshire
    33 = 2805   # New Jersey
    34 = 36     # New Mexico
    35 = 10834  # New York
    36 = 113    # North Carolina
    37 = 9      # North Dakota
    38 = 324    # Ohio
    39 = 108    # Oklahoma
    40 = 55     # Oregon
    41 = 696    # Pennsylvania
    42 = 45     # Puerto Rico
    43 = 80     # Rhode Island
    44 = 97     # South Carolina
    45 = 6      # South Dakota
    46 = 124    # Tennessee
    47 = 350    # Texas
    48 = 19     # Utah
    49 = 29     # Vermont
    50 = 154    # Virginia
    51 = 547    # Washington
    52 = 10     # West Virginia
    53 = 170    # Wisconsin
    54 = 1      # Wyoming
}

foreach ($row in $values.Keys) {
    $ws.Range("AG$row").Value = $values[$row]
}

# Match the look of the new column to the column immediately to its left
# (border / number format / fill), since it's a straight continuation of
# the existing table.
$ws.Range("AF1:AF54").Copy()
$ws.Range("AG1:AG54").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Leave the selection on the newly added header's data cell and scroll the
# window so the new column is visible, matching the author's final view.
$ws.Activate()
$ws.Range("AG2").Select()
$excel.ActiveWindow.ScrollColumn = 20
